$wb = $excel.ActiveWorkbook

# Sheets that need the "想去人数" (F column) counts refreshed: 展览 and 全部类型
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 72
    $ws.Range("F3").Value = 1061
    $ws.Range("F7").Value = 2194
    $ws.Range("F8").Value = 188
    $ws.Range("F10").Value = 1033
}
